$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.772.41"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.321.93"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.80"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.79"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.503"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.89"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.75%  "
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.71"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.122"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.69"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.685.16"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.322.35"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.787"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.706.73"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.99"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.33%  "
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0885"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.83"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.56"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.23"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.52"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.62%  "
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.10"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.33"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "139.64"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -16.09%  "
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("E34").Value = "  -3.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0696"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.33"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.81%  "
$ws.Range("E37").Value = "  -4.50%  "
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.35"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +23.60%  "
$ws.Range("E41").Value = "  -2.64%  "
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.932.65"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.22%  "
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.20"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.70"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.552.40"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.52"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.08"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.20%  "
